$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2010-18")

# New row of data (row 5): CW3M C493 / Demo_Baseline WRB 2010-18 / 2010-18
$ws.Range("A5").Value = "CW3M C493"
$ws.Range("B5").Value = "Demo_Baseline WRB 2010-18"
$ws.Range("C5").Value = "2010-18"

$ws.Range("D5").Value = 1150.4691636666666
$ws.Range("E5").Value = 1612.6987305555554
$ws.Range("F5").Value = 14.557834333333334
$ws.Range("G5").Value = 52.671807666666659
$ws.Range("H5").Value = 5.2565644444444439
$ws.Range("I5").Value = 9.0018087777777769
$ws.Range("J5").Value = 2.782013222222222
$ws.Range("K5").Value = 611.258599
$ws.Range("L5").Value = 44.391417555555549
$ws.Range("M5").Value = 1017.0443863333335
$ws.Range("N5").Value = 1168.4154052222223
$ws.Range("O5").Value = 517773.97222222225
$ws.Range("P5").Value = 286785.73958333331
$ws.Range("Q5").Value = -0.76408811111111108
$ws.Range("R5").Value = -0.00027411111111111109

# Match formatting of the rest of the table (numeric number-formats on D:R,
# and the centered style used in column C).
$ws.Range("D2:R2").Copy()
$ws.Range("D5:R5").PasteSpecial(-4122)
$ws.Range("C4").Copy()
$ws.Range("C5").PasteSpecial(-4122)

$excel.CutCopyMode = 0

$ws.Range("C5").Select()
